$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Text (column B) updates
$ws.Range("B4").Value = "<many>"
$ws.Range("B9").Value = "<your>"
$ws.Range("B18").Value = "<all>"
$ws.Range("B24").Value = "<there>"
$ws.Range("B33").Value = "<line>"
$ws.Range("B48").Value = "<alt>"
$ws.Range("B50").Value = "<xoranwar>"

# Numeric (column C) updates
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 9
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 3
$ws.Range("C12").Value = 3
$ws.Range("C13").Value = 11
$ws.Range("C14").Value = 9
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 5
$ws.Range("C18").Value = 11
$ws.Range("C19").Value = 10
$ws.Range("C20").Value = 8
$ws.Range("C21").Value = 6
$ws.Range("C22").Value = 4
$ws.Range("C23").Value = 6
$ws.Range("C24").Value = 8
$ws.Range("C25").Value = 5
$ws.Range("C26").Value = 8
$ws.Range("C29").Value = 8
$ws.Range("C30").Value = 16
$ws.Range("C32").Value = 10
$ws.Range("C34").Value = 10
$ws.Range("C35").Value = 7
$ws.Range("C36").Value = 10
$ws.Range("C37").Value = 13
$ws.Range("C38").Value = 4
$ws.Range("C39").Value = 7
$ws.Range("C41").Value = 11
$ws.Range("C42").Value = 8
$ws.Range("C43").Value = 4
$ws.Range("C44").Value = 3
$ws.Range("C46").Value = 14
$ws.Range("C47").Value = 8
$ws.Range("C49").Value = 4
$ws.Range("C50").Value = 9
$ws.Range("C51").Value = 7
$ws.Range("C52").Value = 4
